$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "329.82"
$cell.ClearFormats()
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "6.36%"
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "40.15"
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "6.94%"
$cell.ClearFormats()
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "5.274"
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "1.72%"
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "0.08097"
$cell.ClearFormats()
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "2.86%"
$cell.ClearFormats()
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "4.527"
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "2.43%"
$cell.ClearFormats()
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "8.650"
$cell.ClearFormats()
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "4.69%"
$cell.ClearFormats()
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "1.927"
$cell.ClearFormats()
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "0.83%"
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "-1.44%"
$cell.ClearFormats()
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.9364"
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "-0.04%"
$cell.ClearFormats()
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.1338"
$cell.ClearFormats()
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "19.42%"
$cell.ClearFormats()
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.1967"
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "-0.21%"
$cell.ClearFormats()
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.09093"
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "-0.09%"
$cell.ClearFormats()
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.03501"
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "6.46%"
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.09585"
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "-0.19%"
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.001403"
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "1.79%"
$cell.ClearFormats()
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.006307"
$cell.ClearFormats()
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "3.16%"
$cell.ClearFormats()
$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.004319"
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "-5.68%"
$cell.ClearFormats()
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "3.359"
$cell.ClearFormats()
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "-7.01%"
$cell.ClearFormats()
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.3520"
$cell.ClearFormats()
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "3.19%"
$cell.ClearFormats()
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.500"
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "0.68%"
$cell.ClearFormats()
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.1325"
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "3.37%"
$cell.ClearFormats()
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.2572"
$cell.ClearFormats()
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "2.14%"
$cell.ClearFormats()
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "0.04439"
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "1.11%"
$cell.ClearFormats()
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.001225"
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "-0.78%"
$cell.ClearFormats()
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.0001292"
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "-5.10%"
$cell.ClearFormats()
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.0003996"
$cell.ClearFormats()
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "0.11%"
$cell.ClearFormats()
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.02503"
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "13.20%"
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.05190"
$cell.ClearFormats()
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.007698"
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "3.17%"
$cell.ClearFormats()
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.1426"
$cell.ClearFormats()
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "5.38%"
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.009147"
$cell.ClearFormats()
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "4.47%"
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.002164"
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "1.46%"
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.008991"
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "4.21%"
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.00006631"
$cell.ClearFormats()
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "1.06%"
$cell.ClearFormats()
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000751"
$cell.ClearFormats()
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "0.00%"
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.003346"
$cell.ClearFormats()
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "16.91%"
$cell.ClearFormats()
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "147.84%"
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.00002104"
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "0.00%"
$cell.ClearFormats()
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0002003"
$cell.ClearFormats()
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "0.00%"
$cell.ClearFormats()
